$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L112").Value = 106240
$ws.Range("L113").Value = 208584
$ws.Range("L114").Value = 292044
$ws.Range("L115").Value = 349352
$ws.Range("L116").Value = 352093
$ws.Range("G160").Value = 731835904282
$ws.Range("L160").Value = 731835904282

for ($r = 137; $r -le 160; $r++) {
    $ws.Cells.Item($r, 15).Value = "VERIFICAR FORMA DE ACESSO"
}
for ($r = 161; $r -le 178; $r++) {
    $ws.Cells.Item($r, 15).Value = "VÁLIDO"
}
